$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 5.560500000000003
$ws.Range("B12").Value = 4.992100000000001
$ws.Range("B18").Value = 6.541099999999995
$ws.Range("B37").Value = 8.847100000000005
$ws.Range("B55").Value = 6.575899999999991
$ws.Range("B68").Value = 4.761099999999995
$ws.Range("B77").Value = 9.505900000000008
$ws.Range("B78").Value = 9.787600000000001
